$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atribuicoes")

# Add new column I values to row 2 (Tarefas) and row 3 (Técnicos atribuídos)
$ws.Range("I2").Value = 8
$ws.Range("I2").Borders.LineStyle = 1

# Update existing row 3 values
$ws.Range("B3").Value = 5
$ws.Range("E3").Value = 4
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 1
$ws.Range("I3").Borders.LineStyle = 1
